$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 28 -> 27.85546875 (closest achievable quantum is 27.833333...) ---
$ws.Columns.Item(1).ColumnWidth = 27

# --- New cell style (cellXfs index 6): numFmt #,##0 (builtin 3), font "Aptos Narrow"
#     (same font used by existing style index 5, fontId 3), horizontal=left,
#     vertical=center.
#
#     Build the finished style once on an unused scratch cell (starting from
#     a cell that already uses the "Aptos Narrow" font so that font gets
#     reused rather than a new font entry being created), then copy that
#     finished format onto both A129 and C129. Copying the already-finished
#     format (rather than re-running the property assignments on a second
#     cell) avoids leaving a duplicate/orphan style behind in cellXfs. ---
$scratch = $ws.Range("Z500")
$ws.Range("A46").Copy()
$scratch.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$scratch.NumberFormat = "#,##0"
$scratch.HorizontalAlignment = -4131
$scratch.VerticalAlignment = -4108

$scratch.Copy()
$ws.Range("A129").PasteSpecial(-4122)
$ws.Range("C129").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$scratch.Clear()

# --- Selection moves from C127 to D126 ---
$ws.Range("D126").Select()
